$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RMA Details Maintenance Grid")

# Row 2 (RMA test case 1): Sales Order Line, Shipper Line, Id
$ws.Range("E2").Value = "RMA-CFAT-001"
$ws.Range("F2").Value = "RMA-CFAT-1-1"
$ws.Range("J2").Value = "a7s5f000000xLLHAA2"

# Row 3 (RMA test case 2): Sales Order Line, Shipper Line, Id
$ws.Range("E3").Value = "RMA-CFAT-002"
$ws.Range("F3").Value = "RMA-CFAT-1-2"
$ws.Range("J3").Value = "a7s5f000000xLLIAA2"

# Row 4 (RMA test case 3): Sales Order Line, Shipper Line, Id
$ws.Range("E4").Value = "RMA-CFAT-003"
$ws.Range("F4").Value = "RMA-CFAT-1-3"
$ws.Range("J4").Value = "a7s5f000000xLLJAA2"
